# Horarios Linea 141 - scrape refresh 16:43:37 -> 16:53:01
# Updates the "LP1912" sheet's data rows (re-sorted/refreshed arrival log),
# appends the newly scraped rows, and refreshes the "last updated" / "total
# rows" banner text on all three sheets.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item(1)   # LP1912
$ws2 = $wb.Worksheets.Item(2)   # LP1912-215
$ws3 = $wb.Worksheets.Item(3)   # 6203-6173

# --- Banner rows (row 2 = last-updated timestamp, row 3 = row count) ------
$ws1.Range("A2").Value = "Última actualización: 16:53:01"
$ws1.Range("A3").Value = "Total filas: 385"

$ws2.Range("A2").Value = "Última actualización: 16:53:01"
$ws3.Range("A2").Value = "Última actualización: 16:53:01"

# --- Data rows that changed on the LP1912 sheet ---------------------------
# Each entry is the final (post-refresh) content for that row: Hora_Scrap,
# Hora_Llegada, Linea, Minutos, Parada. Rows 385-390 are brand-new rows
# appended by this refresh; the rest are existing rows whose values shifted
# around because the underlying log got re-sorted/merged with new scrapes.
$rows = @(
    @{ Row = 77;  A = "07:20:40"; B = "08:03"; C = "11_ETCHEVERRY";             D = 43;  E = "LP1912" }
    @{ Row = 78;  A = "06:43:40"; B = "08:03"; C = "23_HERNANDEZ";              D = 80;  E = "LP1912" }
    @{ Row = 209; A = "12:37:14"; B = "12:37"; C = "23_HERNANDEZ";              D = 0;   E = "LP1912" }
    @{ Row = 210; A = "11:20:07"; B = "12:37"; C = "27_EL RETIRO";              D = 77;  E = "LP1912" }
    @{ Row = 245; A = "12:01:50"; B = "13:50"; C = "215A_EL PATO";              D = 109; E = "LP1912" }
    @{ Row = 246; A = "13:19:56"; B = "13:50"; C = "11_ETCHEVERRY";             D = 31;  E = "LP1912" }
    @{ Row = 329; A = "15:57:48"; B = "16:29"; C = "10_OLMOS";                  D = 32;  E = "LP1912" }
    @{ Row = 330; A = "15:31:33"; B = "16:29"; C = "14_ABASTO";                 D = 58;  E = "LP1912" }
    @{ Row = 341; A = "16:33:08"; B = "16:43"; C = "16_P MOR-SANTA ANA";        D = 10;  E = "LP1912" }
    @{ Row = 342; A = "14:58:43"; B = "16:43"; C = "225_GOMEZ";                 D = 105; E = "LP1912" }
    @{ Row = 350; A = "16:18:55"; B = "17:04"; C = "23_HERNANDEZ";              D = 46;  E = "LP1912" }
    @{ Row = 351; A = "15:31:33"; B = "17:04"; C = "215A_EL PATO";              D = 93;  E = "LP1912" }
    @{ Row = 352; A = "15:57:48"; B = "17:04"; C = "11_ETCHEVERRY";             D = 67;  E = "LP1912" }
    @{ Row = 358; A = "16:53:01"; B = "17:17"; C = "11_ETCHEVERRY";             D = 24;  E = "LP1912" }
    @{ Row = 359; A = "16:18:55"; B = "17:20"; C = "16_SANTA ANA";              D = 62;  E = "LP1912" }
    @{ Row = 360; A = "16:18:55"; B = "17:20"; C = "26_HERNANDEZ";              D = 62;  E = "LP1912" }
    @{ Row = 361; A = "15:31:33"; B = "17:21"; C = "26_HERNANDEZ";              D = 110; E = "LP1912" }
    @{ Row = 362; A = "16:33:08"; B = "17:21"; C = "16_SANTA ANA";              D = 48;  E = "LP1912" }
    @{ Row = 363; A = "15:31:33"; B = "17:24"; C = "84_COLONIA URQUIZA-ESC 49"; D = 113; E = "LP1912" }
    @{ Row = 364; A = "16:18:55"; B = "17:28"; C = "14_ABASTO";                 D = 70;  E = "LP1912" }
    @{ Row = 365; A = "16:33:08"; B = "17:29"; C = "14_ABASTO";                 D = 56;  E = "LP1912" }
    @{ Row = 366; A = "16:18:55"; B = "17:30"; C = "27_EL RETIRO";              D = 72;  E = "LP1912" }
    @{ Row = 367; A = "16:33:08"; B = "17:31"; C = "15_ABASTO";                 D = 58;  E = "LP1912" }
    @{ Row = 368; A = "16:33:08"; B = "17:33"; C = "27_EL RETIRO";              D = 60;  E = "LP1912" }
    @{ Row = 369; A = "16:43:37"; B = "17:34"; C = "23_HERNANDEZ";              D = 51;  E = "LP1912" }
    @{ Row = 370; A = "16:53:01"; B = "17:34"; C = "10_OLMOS";                  D = 41;  E = "LP1912" }
    @{ Row = 371; A = "15:57:48"; B = "17:35"; C = "27_EL RETIRO";              D = 98;  E = "LP1912" }
    @{ Row = 372; A = "16:53:01"; B = "17:35"; C = "23_HERNANDEZ";              D = 42;  E = "LP1912" }
    @{ Row = 373; A = "16:43:37"; B = "17:36"; C = "27_EL RETIRO";              D = 53;  E = "LP1912" }
    @{ Row = 374; A = "16:53:01"; B = "17:37"; C = "27_EL RETIRO";              D = 44;  E = "LP1912" }
    @{ Row = 375; A = "16:18:55"; B = "17:38"; C = "17_ROMERO";                 D = 80;  E = "LP1912" }
    @{ Row = 376; A = "16:18:55"; B = "17:39"; C = "215B_EL PATO";              D = 81;  E = "LP1912" }
    @{ Row = 377; A = "15:57:48"; B = "17:40"; C = "215B_EL PATO";              D = 103; E = "LP1912" }
    @{ Row = 378; A = "16:43:37"; B = "17:40"; C = "16_SANTA ANA";              D = 57;  E = "LP1912" }
    @{ Row = 379; A = "16:33:08"; B = "17:41"; C = "16_SANTA ANA";              D = 68;  E = "LP1912" }
    @{ Row = 380; A = "15:57:48"; B = "17:41"; C = "17_ROMERO";                 D = 104; E = "LP1912" }
    @{ Row = 381; A = "16:43:37"; B = "17:45"; C = "15_ABASTO";                 D = 62;  E = "LP1912" }
    @{ Row = 382; A = "15:57:48"; B = "17:50"; C = "16_P MOR-167 Y 521";        D = 113; E = "LP1912" }
    @{ Row = 383; A = "16:33:08"; B = "17:51"; C = "16_P MOR-167 Y 521";        D = 78;  E = "LP1912" }
    @{ Row = 384; A = "15:57:48"; B = "17:52"; C = "81_EL PELIGRO";             D = 115; E = "LP1912" }
    @{ Row = 385; A = "16:33:08"; B = "18:04"; C = "17_ROMERO";                 D = 91;  E = "LP1912" }
    @{ Row = 386; A = "16:53:01"; B = "18:09"; C = "14_ABASTO";                 D = 76;  E = "LP1912" }
    @{ Row = 387; A = "16:33:08"; B = "18:21"; C = "26_HERNANDEZ";              D = 108; E = "LP1912" }
    @{ Row = 388; A = "16:33:08"; B = "18:28"; C = "215C_EL PATO";              D = 115; E = "LP1912" }
    @{ Row = 389; A = "16:43:37"; B = "18:32"; C = "11X44_ETCHEVERRY";          D = 109; E = "LP1912" }
    @{ Row = 390; A = "16:53:01"; B = "18:48"; C = "14X44_ABASTO";              D = 115; E = "LP1912" }
)

foreach ($r in $rows) {
    $ws1.Range("A" + $r.Row).Value = $r.A
    $ws1.Range("B" + $r.Row).Value = $r.B
    $ws1.Range("C" + $r.Row).Value = $r.C
    $ws1.Range("D" + $r.Row).Value = $r.D
    $ws1.Range("E" + $r.Row).Value = $r.E
}
